# Update "Kosten in €" (H2:H8) values.
# These cells are stored as text (the source workbook uses inline strings
# for every cell, including the numeric-looking cost figures), so a leading
# apostrophe is used to force a text entry instead of letting Excel coerce
# the input into a number - this mirrors how the original values were typed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = "'0.64"
$ws.Range("H3").Value = "'0.34"
$ws.Range("H4").Value = "'2.56"
$ws.Range("H5").Value = "'3.76"
$ws.Range("H6").Value = "'0.94"
$ws.Range("H7").Value = "'3.76"
$ws.Range("H8").Value = "'1.88"
